# Update cryptocurrency price/volume data to the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.440.23"
$ws.Range("E2").Value = "  +5.81%  "
$ws.Range("D3").Value = "2.622.67"
$ws.Range("E3").Value = "  +8.35%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.40"
$ws.Range("E5").Value = "  +3.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.63"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  -3.99%  "
$ws.Range("D9").Value = "2.660.03"
$ws.Range("E9").Value = "  +8.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.49"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("E11").Value = "  +5.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.344"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "3.071.88"
$ws.Range("E14").Value = "  +7.80%  "
$ws.Range("D15").Value = "60.519.85"
$ws.Range("E15").Value = "  +5.98%  "
$ws.Range("E16").Value = "  +5.24%  "
$ws.Range("E17").Value = "  +5.93%  "
$ws.Range("D18").Value = "2.652.14"
$ws.Range("E18").Value = "  +8.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.59"
$ws.Range("E20").Value = "  +6.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.52"
$ws.Range("E21").Value = "  +4.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.20"
$ws.Range("E22").Value = "  +4.38%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.06"
$ws.Range("E24").Value = "  +3.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"
$ws.Range("E25").Value = "  +4.84%  "
$ws.Range("E26").Value = "  +3.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.990"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  +9.86%  "
$ws.Range("E29").Value = "  +4.49%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.53"
$ws.Range("E31").Value = "  +4.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.49"
$ws.Range("E32").Value = "  +3.56%  "
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("E34").Value = "  +9.13%  "
$ws.Range("E35").Value = "  +7.38%  "
$ws.Range("E36").Value = "  +4.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "308.90"
$ws.Range("E37").Value = "  +10.03%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").Value = "  +8.67%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.855"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.840"
$ws.Range("E40").Value = "  +29.39%  "
$ws.Range("E41").Value = "  +6.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.41"
$ws.Range("E42").Value = "  +3.88%  "
$ws.Range("E43").Value = "  +5.85%  "
$ws.Range("E44").Value = "  +8.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.100"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("E46").Value = "  +14.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.992"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.94"
$ws.Range("E48").Value = "  +7.51%  "
$ws.Range("E49").Value = "  +4.26%  "
$ws.Range("D50").Value = "2.056.24"
$ws.Range("E50").Value = "  +8.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.28"
$ws.Range("E51").Value = "  +0.81%  "
